# Day 116 to Day 120 "Binary Search Tree" completed:
# Mark column C ("Done [yes or no]") as "yes" for rows 214-235 (the
# "Binary Search Trees" topic block), each keeping the per-row highlight
# color (cell style) the author used elsewhere in the sheet for a "yes" cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells already in the sheet that carry the three distinct "yes" fill
# styles needed (12, 14, 15) - used as format-paint sources so the exact
# same style index is reused rather than a new one being created.
$styleSource = @{
    12 = "C13"
    14 = "C20"
    15 = "C22"
}

# Target style per row, taken from the diff.
$rowStyles = [ordered]@{
    214 = 15
    215 = 12
    216 = 15
    217 = 15
    218 = 15
    219 = 12
    220 = 15
    221 = 12
    222 = 15
    223 = 15
    224 = 14
    225 = 12
    226 = 15
    227 = 12
    228 = 12
    229 = 15
    230 = 15
    231 = 12
    232 = 14
    233 = 15
    234 = 14
    235 = 14
}

foreach ($row in $rowStyles.Keys) {
    $style = $rowStyles[$row]
    $srcRef = $styleSource[$style]
    $dst = $ws.Range("C$row")

    # Copy just the formatting (fill colour / style index) from the
    # matching source cell, then set the value - keeps the same shared
    # string ("yes") that's already used throughout the sheet.
    $ws.Range($srcRef).Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $dst.Value = "yes"
}

$excel.CutCopyMode = 0

# The author's selection ended on C234 after marking the section done.
[void]$ws.Range("C234").Select()
